# Updated cryptos list values (Price / Volume(1h)) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddress, $text) {
    $rng = $ws.Range($cellAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '27.113.01'
Set-TextValue "E2" '  +0.62%  '

# Row 3
Set-TextValue "D3" '1.682.16'
Set-TextValue "E3" '  +0.52%  '

# Row 4
Set-TextValue "E4" '  +0.14%  '

# Row 5
Set-TextValue "D5" '215.25'
Set-TextValue "E5" '  +0.19%  '

# Row 6
Set-TextValue "E6" '  +0.20%  '

# Row 7
Set-TextValue "E7" '  +0.07%  '

# Row 9
Set-TextValue "D9" '21.32'
Set-TextValue "E9" '  +5.61%  '

# Row 10
Set-TextValue "E10" '  +0.66%  '

# Row 11
Set-TextValue "E11" '  -0.47%  '

# Row 12
Set-TextValue "D12" '1.918.84'
Set-TextValue "E12" '  +0.51%  '

# Row 13
Set-TextValue "D13" '1.692.26'
Set-TextValue "E13" '  +1.20%  '

# Row 14
Set-TextValue "E14" '  +1.61%  '

# Row 15
Set-TextValue "E15" '  +2.07%  '

# Row 16
Set-TextValue "D16" '66.22'
Set-TextValue "E16" '  +0.83%  '

# Row 17
Set-TextValue "D17" '27.106.08'
Set-TextValue "E17" '  +0.59%  '

# Row 18
Set-TextValue "E18" '  +1.64%  '

# Row 19
Set-TextValue "D19" '8.12'
Set-TextValue "E19" '  +0.38%  '

# Row 20
Set-TextValue "D20" '0.0₃0750'
Set-TextValue "E20" '  +2.30%  '

# Row 21
Set-TextValue "E21" '  +0.11%  '

# Row 22
Set-TextValue "D22" '4.51'
Set-TextValue "E22" '  +1.61%  '

# Row 23
Set-TextValue "E23" '  +2.57%  '

# Row 24
Set-TextValue "E24" '  -2.52%  '

# Row 25
Set-TextValue "D25" '146.91'
Set-TextValue "E25" '  +0.93%  '

# Row 26
Set-TextValue "D26" '7.22'
Set-TextValue "E26" '  +0.93%  '

# Row 27
Set-TextValue "D27" '16.32'
Set-TextValue "E27" '  +2.14%  '

# Row 28
Set-TextValue "E28" '  +0.67%  '

# Row 29
Set-TextValue "E29" '  +0.14%  '

# Row 30
Set-TextValue "E30" '  +0.48%  '

# Row 31
Set-TextValue "E31" '  +0.11%  '

# Row 32
Set-TextValue "D32" '1.559.24'
Set-TextValue "E32" '  +5.62%  '

# Row 33
Set-TextValue "E33" '  +1.01%  '

# Row 34
Set-TextValue "D34" '3.19'
Set-TextValue "E34" '  +1.76%  '

# Row 35
Set-TextValue "E35" '  +2.43%  '

# Row 36
Set-TextValue "E36" '  +4.49%  '

# Row 37
Set-TextValue "D37" '0.938'
Set-TextValue "E37" '  +4.68%  '

# Row 38
Set-TextValue "E38" '  -1.26%  '

# Row 39
Set-TextValue "E39" '  +2.44%  '

# Row 40
Set-TextValue "E40" '  -0.37%  '

# Row 41
Set-TextValue "D41" '68.87'
Set-TextValue "E41" '  +3.23%  '

# Row 42
Set-TextValue "E42" '  +0.12%  '

# Row 43
Set-TextValue "E43" '  -2.81%  '

# Row 44
Set-TextValue "E44" '  -1.59%  '

# Row 45
Set-TextValue "D45" '1.827.27'
Set-TextValue "E45" '  +0.74%  '

# Row 46
Set-TextValue "E46" '  +0.39%  '

# Row 47
Set-TextValue "D47" '90.81'
Set-TextValue "E47" '  +0.43%  '

# Row 48
Set-TextValue "E48" '  +3.16%  '

# Row 49
Set-TextValue "D49" '0.0₆0108'
Set-TextValue "E49" '  +1.59%  '

# Row 50
Set-TextValue "E50" '  +3.51%  '

# Row 51
Set-TextValue "D51" '8.07'
Set-TextValue "E51" '  +4.95%  '
